$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current row 741 (shifts 741:776 down to 745:780)
$ws.Range("741:744").Insert()

# Row 741 (new): Lapins / Especial
$ws.Range("A741").Value = 6
$ws.Range("B741").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C741").Value = "Metropolitana"
$ws.Range("D741").Value = 44578
$ws.Range("E741").Value = 13
$ws.Range("F741").Value = "Fruta"
$ws.Range("G741").Value = 100103
$ws.Range("H741").Value = "Frutos de hueso (carozo)"
$ws.Range("I741").Value = 100103001
$ws.Range("J741").Value = "Cereza"
$ws.Range("K741").Value = "Lapins"
$ws.Range("L741").Value = "Especial"
$ws.Range("M741").Value = 300
$ws.Range("N741").Value = 9000
$ws.Range("O741").Value = 9000
$ws.Range("P741").Value = 9000
$ws.Range("Q741").Value = "$/caja 15 kilos"
$ws.Range("R741").Value = "Provincia de Curicó"
$ws.Range("S741").Value = 600
$ws.Range("T741").Value = 15

# Row 742 (new): Santina / Segunda
$ws.Range("A742").Value = 6
$ws.Range("B742").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C742").Value = "Metropolitana"
$ws.Range("D742").Value = 44578
$ws.Range("E742").Value = 13
$ws.Range("F742").Value = "Fruta"
$ws.Range("G742").Value = 100103
$ws.Range("H742").Value = "Frutos de hueso (carozo)"
$ws.Range("I742").Value = 100103001
$ws.Range("J742").Value = "Cereza"
$ws.Range("K742").Value = "Santina"
$ws.Range("L742").Value = "Segunda"
$ws.Range("M742").Value = 480
$ws.Range("N742").Value = 3000
$ws.Range("O742").Value = 4000
$ws.Range("P742").Value = 3500
$ws.Range("Q742").Value = "$/bandeja 10 kilos"
$ws.Range("R742").Value = "Provincia de Curicó"
$ws.Range("S742").Value = 350
$ws.Range("T742").Value = 10

# Row 743 (new): Santina / Tercera
$ws.Range("A743").Value = 6
$ws.Range("B743").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C743").Value = "Metropolitana"
$ws.Range("D743").Value = 44578
$ws.Range("E743").Value = 13
$ws.Range("F743").Value = "Fruta"
$ws.Range("G743").Value = 100103
$ws.Range("H743").Value = "Frutos de hueso (carozo)"
$ws.Range("I743").Value = 100103001
$ws.Range("J743").Value = "Cereza"
$ws.Range("K743").Value = "Santina"
$ws.Range("L743").Value = "Tercera"
$ws.Range("M743").Value = 360
$ws.Range("N743").Value = 2000
$ws.Range("O743").Value = 2000
$ws.Range("P743").Value = 2000
$ws.Range("Q743").Value = "$/bandeja 10 kilos"
$ws.Range("R743").Value = "Provincia de Curicó"
$ws.Range("S743").Value = 200
$ws.Range("T743").Value = 10

# Row 744 (new): Sweet Heart / Primera
$ws.Range("A744").Value = 6
$ws.Range("B744").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C744").Value = "Metropolitana"
$ws.Range("D744").Value = 44578
$ws.Range("E744").Value = 13
$ws.Range("F744").Value = "Fruta"
$ws.Range("G744").Value = 100103
$ws.Range("H744").Value = "Frutos de hueso (carozo)"
$ws.Range("I744").Value = 100103001
$ws.Range("J744").Value = "Cereza"
$ws.Range("K744").Value = "Sweet Heart"
$ws.Range("L744").Value = "Primera"
$ws.Range("M744").Value = 720
$ws.Range("N744").Value = 4000
$ws.Range("O744").Value = 5000
$ws.Range("P744").Value = 4500
$ws.Range("Q744").Value = "$/bandeja 10 kilos"
$ws.Range("R744").Value = "Provincia de Curicó"
$ws.Range("S744").Value = 450
$ws.Range("T744").Value = 10
